$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.165938864628821
$ws.Range("C2").Value = 0.5938864628820961
$ws.Range("J2").Value = 0.008733624454148471
$ws.Range("P2").Value = 0.1091703056768559
$ws.Range("S2").Value = 0.1222707423580786
$ws.Range("B3").Value = 0.01398601398601399
$ws.Range("C3").Value = 0.03496503496503497
$ws.Range("J3").Value = 0.02797202797202797
$ws.Range("P3").Value = 0.7412587412587412
$ws.Range("S3").Value = 0.1818181818181818
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.7027027027027027
$ws.Range("S4").Value = 0.2702702702702703
$ws.Range("B6").Value = 0.05851063829787234
$ws.Range("D6").Value = 0.02659574468085106
$ws.Range("F6").Value = 0.06382978723404255
$ws.Range("J6").Value = 0.25
$ws.Range("O6").Value = 0.02659574468085106
$ws.Range("Q6").Value = 0.1968085106382979
$ws.Range("R6").Value = 0.05851063829787234
$ws.Range("S6").Value = 0.3191489361702128
$ws.Range("B7").Value = 0.1197916666666667
$ws.Range("D7").Value = 0.005208333333333333
$ws.Range("F7").Value = 0.05208333333333334
$ws.Range("J7").Value = 0.171875
$ws.Range("O7").Value = 0.03125
$ws.Range("Q7").Value = 0.1614583333333333
$ws.Range("R7").Value = 0.09375
$ws.Range("S7").Value = 0.3645833333333333
$ws.Range("B8").Value = 0.09713024282560706
$ws.Range("D8").Value = 0.01324503311258278
$ws.Range("F8").Value = 0.05739514348785872
$ws.Range("J8").Value = 0.1103752759381898
$ws.Range("O8").Value = 0.008830022075055188
$ws.Range("Q8").Value = 0.2185430463576159
$ws.Range("R8").Value = 0.08830022075055188
$ws.Range("S8").Value = 0.4061810154525387
$ws.Range("B9").Value = 0.0576923076923077
$ws.Range("F9").Value = 0.0673076923076923
$ws.Range("J9").Value = 0.1153846153846154
$ws.Range("O9").Value = 0.01442307692307692
$ws.Range("Q9").Value = 0.2355769230769231
$ws.Range("R9").Value = 0.07211538461538461
$ws.Range("S9").Value = 0.4375
$ws.Range("B10").Value = 0.08691674290942361
$ws.Range("D10").Value = 0.02195791399817017
$ws.Range("E10").Value = 0.002744739249771272
$ws.Range("F10").Value = 0.0677035681610247
$ws.Range("J10").Value = 0.1189387008234218
$ws.Range("O10").Value = 0.010064043915828
$ws.Range("Q10").Value = 0.2049405306495883
$ws.Range("R10").Value = 0.08600182982616651
$ws.Range("S10").Value = 0.4007319304666057
$ws.Range("G11").Value = 0.1254355400696864
$ws.Range("J11").Value = 0.1045296167247387
$ws.Range("K11").Value = 0.2020905923344948
$ws.Range("L11").Value = 0.554006968641115
$ws.Range("S11").Value = 0.01393728222996516
$ws.Range("G12").Value = 0.7345679012345679
$ws.Range("J12").Value = 0.1666666666666667
$ws.Range("L12").Value = 0.02469135802469136
$ws.Range("S12").Value = 0.07407407407407407
$ws.Range("G13").Value = 0.6984126984126984
$ws.Range("J13").Value = 0.1746031746031746
$ws.Range("S13").Value = 0.126984126984127
$ws.Range("G14").Value = 0.5
$ws.Range("S14").Value = 0.5
$ws.Range("F15").Value = 0.02116402116402116
$ws.Range("H15").Value = 0.1746031746031746
$ws.Range("I15").Value = 0.1005291005291005
$ws.Range("J15").Value = 0.3174603174603174
$ws.Range("K15").Value = 0.07936507936507936
$ws.Range("M15").Value = 0.02645502645502645
$ws.Range("O15").Value = 0.03703703703703703
$ws.Range("S15").Value = 0.2433862433862434
$ws.Range("F16").Value = 0.01298701298701299
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.07792207792207792
$ws.Range("J16").Value = 0.461038961038961
$ws.Range("K16").Value = 0.07792207792207792
$ws.Range("M16").Value = 0.02597402597402598
$ws.Range("N16").Value = 0.006493506493506494
$ws.Range("O16").Value = 0.03896103896103896
$ws.Range("S16").Value = 0.1168831168831169
$ws.Range("F17").Value = 0.011441647597254
$ws.Range("H17").Value = 0.2196796338672769
$ws.Range("I17").Value = 0.09610983981693363
$ws.Range("J17").Value = 0.356979405034325
$ws.Range("K17").Value = 0.09610983981693363
$ws.Range("M17").Value = 0.03661327231121281
$ws.Range("N17").Value = 0.002288329519450801
$ws.Range("O17").Value = 0.05263157894736842
$ws.Range("S17").Value = 0.1281464530892449
$ws.Range("F18").Value = 0.01685393258426966
$ws.Range("H18").Value = 0.2134831460674157
$ws.Range("I18").Value = 0.07303370786516854
$ws.Range("J18").Value = 0.3651685393258427
$ws.Range("K18").Value = 0.1067415730337079
$ws.Range("M18").Value = 0.03932584269662921
$ws.Range("O18").Value = 0.02247191011235955
$ws.Range("S18").Value = 0.1629213483146068
$ws.Range("F19").Value = 0.0129136400322841
$ws.Range("H19").Value = 0.2098466505246166
$ws.Range("I19").Value = 0.1008878127522195
$ws.Range("J19").Value = 0.3284907183212268
$ws.Range("K19").Value = 0.1097659402744148
$ws.Range("M19").Value = 0.02502017756255044
$ws.Range("N19").Value = 0.0008071025020177562
$ws.Range("O19").Value = 0.07667473769168684
$ws.Range("S19").Value = 0.1355932203389831
